# Add a new price-snapshot column "GR" (timestamp 2026-02-06 07:39:23),
# shifting the existing "nom" (GR->GS) and "url_produit" (GS->GT) columns
# one place to the right. The newly inserted column is seeded with the
# last known price (the value already present in the previous snapshot
# column, GQ) for every product row, just like the scraper that produces
# this workbook carries forward unchanged prices.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before GR; this shifts GR->GS and GS->GT and keeps
# the header/data styles of the surrounding cells.
$ws.Columns("GR:GR").Insert()

# New snapshot timestamp header.
$ws.Range("GR1").Value = "2026-02-06 07:39:23"

# Carry the previous snapshot's price (column GQ) forward into the new
# column GR for every data row (2-210). Rows whose latest known price is
# blank stay blank.
$src = $ws.Range("GQ2:GQ210")
$dst = $ws.Range("GR2:GR210")
$dst.Value = $src.Value()
